$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The workbook already has a table named "Table5" on the CustList
# sheet. Excel's ListObjects.Add() auto-assigns the next free
# "TableN" name to a brand new table, which would collide with that
# existing name. Temporarily rename it out of the way so the new
# table we are about to create does not clobber it, then restore it.
# ------------------------------------------------------------------
$custListTable = $wb.Worksheets.Item("CustList").ListObjects.Item(1)
$custListTable.Name = "TempRenameTable5"

# ------------------------------------------------------------------
# Fill in the flight row data on the ActiveFlights sheet.
# ------------------------------------------------------------------
$ws.Range("A1").Value = "FlightId"
$ws.Range("B1").Value = "DepartingFrom"
$ws.Range("C1").Value = "ArrivingAt"
$ws.Range("D1").Value = "DateTime"

$ws.Range("B2").Value = "houston airport"
$ws.Range("C2").Value = "Nebraska airport"

# "123" and "4/6/2023" look like a number/date, so a direct .Value
# assignment would be auto-converted. Build them as text formulas on
# a scratch cell and paste the computed (text) values into place so
# they land as plain shared-string text with no special formatting.
$ws.Range("Z1").Formula = "=""4/6/2023"""
$ws.Range("Z1").Copy()
$ws.Range("D2").PasteSpecial(-4163)

$ws.Range("Z1").Formula = "=""123"""
$ws.Range("Z1").Copy()
$ws.Range("A2").PasteSpecial(-4163)

$ws.Range("Z1").Clear()

# ------------------------------------------------------------------
# Turn the A1:D2 range into a proper Excel table ("Table6").
# ------------------------------------------------------------------
$newTable = $ws.ListObjects.Add(1, $ws.Range("A1:D2"), [System.Type]::Missing, 1)
$newTable.Name = "Table6"

# Restore the original table's name.
$custListTable.Name = "Table5"

# ------------------------------------------------------------------
# Match the column widths / selection recorded for the sheet.
# ------------------------------------------------------------------
$ws.Columns("A:E").ColumnWidth = 10.17

$ws.Range("D11").Select()
